$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing TEXT interpretation (matches the
# source data's inline-string cells) and then strip the temporary '@' number
# format so the cell's style index is left exactly as it was (style 0/General).
function Set-TextValue([object]$range, [string]$value) {
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '69.485.18'
$ws.Range("E2").Value = '  +0.24%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.671.36'
$ws.Range("E3").Value = '  -0.42%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '648.23'
$ws.Range("E5").Value = '  -4.82%  '

# Row 6
Set-TextValue $ws.Range("D6") '159.67'
$ws.Range("E6").Value = '  +0.16%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("E8").Value = '  +0.41%  '

# Row 9
$ws.Range("E9").Value = '  -0.75%  '

# Row 10
Set-TextValue $ws.Range("D10") '7.10'
$ws.Range("E10").Value = '  -0.18%  '

# Row 11
$ws.Range("E11").Value = '  +0.02%  '

# Row 12
$ws.Range("E12").Value = '  -0.16%  '

# Row 13
Set-TextValue $ws.Range("D13") '4.292.01'
$ws.Range("E13").Value = '  -0.39%  '

# Row 14
Set-TextValue $ws.Range("D14") '32.57'
$ws.Range("E14").Value = '  +0.25%  '

# Row 15
Set-TextValue $ws.Range("D15") '3.661.11'
$ws.Range("E15").Value = '  -0.07%  '

# Row 16
Set-TextValue $ws.Range("D16") '69.446.86'
$ws.Range("E16").Value = '  +0.20%  '

# Row 17
$ws.Range("E17").Value = '  +0.84%  '

# Row 18
Set-TextValue $ws.Range("D18") '15.98'
$ws.Range("E18").Value = '  -0.67%  '

# Row 19
Set-TextValue $ws.Range("D19") '6.44'
$ws.Range("E19").Value = '  -0.34%  '

# Row 20
Set-TextValue $ws.Range("D20") '464.54'
$ws.Range("E20").Value = '  -0.87%  '

# Row 21
Set-TextValue $ws.Range("D21") '9.75'
$ws.Range("E21").Value = '  -1.81%  '

# Row 23
$ws.Range("E23").Value = '  -0.49%  '

# Row 24
Set-TextValue $ws.Range("D24") '3.818.94'

# Row 25
$ws.Range("E25").Value = '  -0.01%  '

# Row 27
Set-TextValue $ws.Range("D27") '10.77'
$ws.Range("E27").Value = '  -1.21%  '

# Row 28
Set-TextValue $ws.Range("D28") '8.91'
$ws.Range("E28").Value = '  -2.36%  '

# Row 29
$ws.Range("E29").Value = '  -2.79%  '

# Row 30
$ws.Range("E30").Value = '  -4.09%  '

# Row 31
$ws.Range("E31").Value = '  -0.27%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.995'
$ws.Range("E32").Value = '  -0.45%  '

# Row 33
$ws.Range("E33").Value = '  -2.76%  '

# Row 34
$ws.Range("E34").Value = '  -1.24%  '

# Row 35
$ws.Range("E35").Value = '  +3.59%  '

# Row 36
Set-TextValue $ws.Range("D36") '3.661.93'

# Row 37
Set-TextValue $ws.Range("D37") '8.36'
$ws.Range("E37").Value = '  +0.96%  '

# Row 39
$ws.Range("E39").Value = '  -5.84%  '

# Row 40
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D40") '177.82'
$ws.Range("E40").Value = '  +4.33%  '

# Row 41
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D41") '1.00'
$ws.Range("E41").Value = '  -0.05%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.0896'
$ws.Range("E42").Value = '  -1.07%  '

# Row 43
$ws.Range("E43").Value = '  -3.87%  '

# Row 44
$ws.Range("E44").Value = '  -1.77%  '

# Row 45
Set-TextValue $ws.Range("D45") '46.60'
$ws.Range("E45").Value = '  -2.10%  '

# Row 46
$ws.Range("E46").Value = '  +0.63%  '

# Row 47
$ws.Range("E47").Value = '  -3.22%  '

# Row 48
Set-TextValue $ws.Range("D48") '26.97'
$ws.Range("E48").Value = '  -5.33%  '

# Row 49
$ws.Range("E49").Value = '  -3.42%  '

# Row 50
Set-TextValue $ws.Range("D50") '7.81'
$ws.Range("E50").Value = '  +0.33%  '

# Row 51
$ws.Range("E51").Value = '  -5.56%  '
